$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B11").Value = 5.852
$ws.Range("D11").Value = -7.549999999999999
$ws.Range("B12").Value = 4.951000000000001
$ws.Range("B15").Value = 5.090000000000001
$ws.Range("D23").Value = -8.422000000000001
$ws.Range("B27").Value = 5.424
$ws.Range("B28").Value = 5.88
$ws.Range("D28").Value = -8.198
$ws.Range("B31").Value = 6.032000000000001
$ws.Range("B32").Value = 6.063
$ws.Range("D32").Value = -7.689
$ws.Range("D34").Value = -8.333
$ws.Range("B36").Value = 8.431000000000001
$ws.Range("D36").Value = -7.513000000000001
$ws.Range("D37").Value = -8.190000000000001
$ws.Range("B38").Value = 6.247000000000001
$ws.Range("D42").Value = -8.286999999999999
$ws.Range("B46").Value = 6.259
$ws.Range("D49").Value = -8.330000000000002
$ws.Range("B54").Value = 4.796
$ws.Range("D54").Value = -8.168000000000001
$ws.Range("B55").Value = 4.722
$ws.Range("B56").Value = 4.654000000000001
$ws.Range("B67").Value = 5.516
$ws.Range("B69").Value = 5.366999999999999
$ws.Range("B72").Value = 5.792999999999999
$ws.Range("B73").Value = 6.444000000000001
$ws.Range("D78").Value = -8.382
$ws.Range("D80").Value = -8.134
$ws.Range("B83").Value = 5.831
$ws.Range("B86").Value = 5.02
$ws.Range("B91").Value = 5.667000000000001
$ws.Range("B93").Value = 5.424000000000001
$ws.Range("D97").Value = -7.106
$ws.Range("B99").Value = 4.76
$ws.Range("D99").Value = -8.164
$ws.Range("D100").Value = -7.962000000000001
$ws.Range("D101").Value = -7.776999999999999
$ws.Range("B104").Value = 7.459000000000001
$ws.Range("B105").Value = 8.617000000000001
